# "Generate Report for Handback"
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoffs have now been handed back / synced with en-US:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview sheet + each language sheet)
#   - Latest Target File / Latest Handback File columns are now populated
#   - Latest Handback DateTime is stamped with the actual handback time

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn ----
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

$zhcn.Range("E2").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.md"
$zhcn.Range("F2").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.zh-cn.xlf"
$zhcn.Range("E3").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.md"
$zhcn.Range("F3").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/d3027a9809954156586cc6f868407073f904742f/e2e/7ed9af88-d712-41cc-a903-662a24ade0ac.md", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/553cbae2377e052a21eafd752493a0d0fa0e02f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.zh-cn.xlf", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/d3027a9809954156586cc6f868407073f904742f/e2e/7ed9af88-d712-41cc-a903-662a24ade0ac.md", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/553cbae2377e052a21eafd752493a0d0fa0e02f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.zh-cn.xlf", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.zh-cn.xlf")

$zhcn.Range("G2").Value = "2016-01-27 03:01:56"
$zhcn.Range("G3").Value = "2016-01-27 03:01:56"

# ---- de-de ----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

$dede.Range("E2").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.md"
$dede.Range("F2").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.de-de.xlf"
$dede.Range("E3").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.md"
$dede.Range("F3").Value = "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/d3027a9809954156586cc6f868407073f904742f/e2e/7ed9af88-d712-41cc-a903-662a24ade0ac.md", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a47341907c1fbaa7ef6cc414ce908359e15d7d4b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.de-de.xlf", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/d3027a9809954156586cc6f868407073f904742f/e2e/7ed9af88-d712-41cc-a903-662a24ade0ac.md", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a47341907c1fbaa7ef6cc414ce908359e15d7d4b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.de-de.xlf", "", "", "7ed9af88-d712-41cc-a903-662a24ade0ac.db3888d5b5fe2e0689b9b94e9c607e06d3ef2c65.de-de.xlf")

$dede.Range("G2").Value = "2016-01-27 03:02:20"
$dede.Range("G3").Value = "2016-01-27 03:02:20"
